$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve these cells as literal text (not auto-converted to numbers) by
# forcing a Text number format before the write, then clearing the explicit
# format afterwards so the cell style index is left exactly as it was before.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "58.755.90"
Set-TextValue "E2" "  -0.13%  "
Set-TextValue "D3" "2.304.95"
Set-TextValue "E3" "  -0.14%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "538.27"
Set-TextValue "E5" "  -1.82%  "
Set-TextValue "D6" "131.93"
Set-TextValue "E6" "  +0.15%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "0.588"
Set-TextValue "E8" "  +2.56%  "
Set-TextValue "D9" "2.303.67"
Set-TextValue "E9" "  -0.18%  "
Set-TextValue "E10" "  -1.78%  "
Set-TextValue "D11" "5.48"
Set-TextValue "E11" "  -1.47%  "
Set-TextValue "E12" "  +0.78%  "
Set-TextValue "D13" "0.333"
Set-TextValue "E13" "  -0.52%  "
Set-TextValue "E14" "  -1.26%  "
Set-TextValue "D15" "2.717.77"
Set-TextValue "E15" "  -0.35%  "
Set-TextValue "D16" "58.648.98"
Set-TextValue "E16" "  -0.28%  "
Set-TextValue "E17" "  -0.37%  "
Set-TextValue "D18" "2.304.97"
Set-TextValue "E18" "  +3.15%  "
Set-TextValue "D19" "10.61"
Set-TextValue "E19" "  -1.14%  "
Set-TextValue "D20" "4.17"
Set-TextValue "E20" "  -3.46%  "
Set-TextValue "D21" "313.60"
Set-TextValue "E21" "  -0.50%  "
Set-TextValue "D22" "6.63"
Set-TextValue "E22" "  +2.34%  "
Set-TextValue "E23" "  +0.04%  "
Set-TextValue "D24" "62.53"
Set-TextValue "E24" "  -1.40%  "
Set-TextValue "D25" "0.173"
Set-TextValue "E25" "  +1.02%  "
Set-TextValue "D27" "7.91"
Set-TextValue "E27" "  -2.53%  "
Set-TextValue "D28" "1.29"
Set-TextValue "E28" "  -1.72%  "
Set-TextValue "D29" "171.10"
Set-TextValue "E29" "  +1.02%  "
Set-TextValue "E30" "  -2.27%  "
Set-TextValue "D31" "0.0₃0732"
Set-TextValue "E31" "  +0.44%  "
Set-TextValue "E32" "  +2.40%  "
Set-TextValue "D33" "5.85"
Set-TextValue "E33" "  +0.83%  "
Set-TextValue "E34" "  +0.49%  "
Set-TextValue "B35" "USDe"
Set-TextValue "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  -0.01%  "
Set-TextValue "B36" "EthereumClassic"
Set-TextValue "C36" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D36" "17.91"
Set-TextValue "E36" "  +0.73%  "
Set-TextValue "D37" "1.31"
Set-TextValue "E37" "  +3.72%  "
Set-TextValue "E38" "  -0.04%  "
Set-TextValue "D39" "4.05"
Set-TextValue "E39" "  +1.66%  "
Set-TextValue "D40" "1.51"
Set-TextValue "E40" "  +0.04%  "
Set-TextValue "D41" "296.53"
Set-TextValue "E41" "  -1.61%  "
Set-TextValue "D42" "141.42"
Set-TextValue "E42" "  -0.20%  "
Set-TextValue "D43" "3.44"
Set-TextValue "E43" "  +0.06%  "
Set-TextValue "E44" "  +0.64%  "
Set-TextValue "D45" "0.0495"
Set-TextValue "E45" "  -1.48%  "
Set-TextValue "E46" "  -0.31%  "
Set-TextValue "D47" "18.26"
Set-TextValue "E47" "  -1.56%  "
Set-TextValue "E48" "  -2.33%  "
Set-TextValue "E49" "  -0.08%  "
Set-TextValue "E50" "  -0.04%  "
Set-TextValue "E51" "  +0.66%  "
